# Update with restock suggestion
# Applies per-week recalculated Inventory Coverage / Seasonality Index / Stockout Risk /
# Reorder Urgency / Lifecycle Stage values to the "Forecast Comparison" sheet, fills in
# the Week_Start_Date column, drops the "Sales Volume Rank" column, and marks the
# Max/Min Forecast Week metrics on the "Summary" sheet as not applicable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# ---------------------------------------------------------------------------
# Week_Start_Date (column B) for rows 2-17 (W1-W16)
# ---------------------------------------------------------------------------
$weekStartDates = @(
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27",
    "2025-05-04",
    "2025-05-11",
    "2025-05-18"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    # Force the value to be stored as literal text instead of letting Excel
    # auto-convert date-shaped strings (e.g. "2025-02-02") into date serials.
    $cell.NumberFormat = "@"
    $cell.Value = $weekStartDates[$i]
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# Inventory Coverage (column L) for rows 2-17
# ---------------------------------------------------------------------------
$inventoryCoverage = @(8.35, 6.56, 4.62, 3.49, 2.46, 1.49, 0.47, 0, 0, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $inventoryCoverage.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $inventoryCoverage[$i]
}

# ---------------------------------------------------------------------------
# Stockout Risk (column M) / Reorder Urgency (column N) for rows 8-17 (W7-W16)
# ---------------------------------------------------------------------------
for ($row = 8; $row -le 17; $row++) {
    $ws.Cells.Item($row, 13).Value = "High"
    $ws.Cells.Item($row, 14).Value = "Urgent"
}

# ---------------------------------------------------------------------------
# Seasonality Index (column P) for rows 2-17 (row 5 and row 10 unchanged)
# ---------------------------------------------------------------------------
$seasonalityIndex = @{
    2  = 1.16
    3  = 0.87
    4  = 0.95
    6  = 1.13
    7  = 0.8100000000000001
    8  = 0.87
    9  = 0.8
    11 = 0.83
    12 = 1.07
    13 = 0.84
    14 = 1.13
    15 = 1.02
    16 = 0.84
    17 = 1.03
}

foreach ($row in $seasonalityIndex.Keys) {
    $ws.Cells.Item($row, 16).Value = $seasonalityIndex[$row]
}

# ---------------------------------------------------------------------------
# Remove the "Sales Volume Rank" column (Q). This shifts the old "Lifecycle
# Stage" column (R) left into Q.
# ---------------------------------------------------------------------------
$ws.Columns.Item(17).Delete()

# ---------------------------------------------------------------------------
# Lifecycle Stage (now column Q after the delete above) for rows 2-17:
# all weeks move from "Growth" to "Decline".
# ---------------------------------------------------------------------------
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 17).Value = "Decline"
}

# ---------------------------------------------------------------------------
# Summary sheet: Max/Min Forecast Week become "N/A"
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B13").Value = "N/A"
$summary.Range("B15").Value = "N/A"
